# Diretivas.docx edit script
# 1. Bold the first "orderBy" word in the "orderBy utiliza com filter..." paragraph.
# 2. Add three new glossary-style paragraphs (currency, number, limitTo) after the
#    "orderBy ... nome'" paragraph, followed by a blank paragraph and a paragraph
#    that now hosts the (previously inline) _GoBack bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: bold the first occurrence of "orderBy" inside the target paragraph.
# ---------------------------------------------------------------------------
$orderByPara = $d.Paragraphs.Item(29)
$findRange = $orderByPara.Range.Duplicate
$findRange.Find.Execute("orderBy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Bold = 1

# ---------------------------------------------------------------------------
# Step 1b: the inline _GoBack bookmark currently sits between "+" and "nome'"
# near the end of that same paragraph; it is relocated to a paragraph of its
# own further down, so strip it out here (re-writing the plain "+nome'" text
# without the bookmark wrapper).
# ---------------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item(29)
$bookmarkFind = $bookmarkPara.Range.Duplicate
$bookmarkFind.Find.Execute("+nome'", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$plusNomeRange = $d.Range($bookmarkFind.Start, $bookmarkFind.End)
$plusNomeXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>+</w:t></w:r><w:r><w:t>nome&apos;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$plusNomeRange.InsertXML($plusNomeXml)

# ---------------------------------------------------------------------------
# Step 2: rebuild the three trailing empty paragraphs into glossary entries,
# a blank paragraph, and a paragraph that carries the _GoBack bookmark.
# ---------------------------------------------------------------------------
$targetPara = $d.Paragraphs.Item(30)
$targetRange = $targetPara.Range

$newBody = @'
<w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>currency</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:r><w:t>refere-se a padronizar estilo de moedas</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>number</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> : {100.26 : number:1} -&gt; imprime com uma casa decimal: 100.2</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>limitTo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">: limita a apresenta&#231;&#227;o de nomes de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arrays</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> na tela</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$xmlPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $newBody + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xmlPackage)

# The two other originally-empty underline paragraphs are no longer needed;
# they now sit right after the five paragraphs we just inserted (index 35/36).
$d.Paragraphs.Item(35).Range.Delete()
$d.Paragraphs.Item(35).Range.Delete()
